$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record needs to be inserted as the first entry of this
# block (row 38), pushing the existing rows 38-44 down to 39-45.
$ws.Rows("38:38").Insert()

# Populate the newly inserted row 38 with the latest "Primera" quality record
# (same market/region/category as the rest of the block, new date & prices).
$ws.Range("A38").Value = 11
$ws.Range("B38").Value = "Vega Monumental Concepción"
$ws.Range("C38").Value = "Bíobío"
$ws.Range("D38").Value = 44694
$ws.Range("E38").Value = 8
$ws.Range("F38").Value = 100112043
$ws.Range("G38").Value = "Pepino dulce"
$ws.Range("H38").Value = "Cultivar IV Región"
$ws.Range("I38").Value = "Primera"
$ws.Range("J38").Value = 100
$ws.Range("K38").Value = 13000
$ws.Range("L38").Value = 14000
$ws.Range("M38").Value = 13500
$ws.Range("N38").Value = "$/bandeja 18 kilos"
$ws.Range("O38").Value = "Provincia de Limarí"
$ws.Range("P38").Value = 750
$ws.Range("Q38").Value = 18
$ws.Range("R38").Value = "Hortaliza"
